$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink attached to the last row's teams-chat cell (I103) before
# the row shift, since Excel's row delete does not relocate hyperlink anchors
# automatically in this engine.
$oldLink = $ws.Range("I103")
$oldLinkAddress = $null
if ($oldLink.Hyperlinks.Count -gt 0) {
    $oldLinkAddress = $oldLink.Hyperlinks.Item(1).Address
    $oldLink.Hyperlinks.Delete()
}

# Delete the "Poste vacant" row (row 101): Nom=vacant, Prenom=Poste, ...
$ws.Rows.Item(101).Delete()

# Re-anchor the hyperlink to its new location (row shifted up by one: I103 -> I102)
if ($oldLinkAddress) {
    $ws.Hyperlinks.Add($ws.Range("I102"), $oldLinkAddress) | Out-Null
}
